$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value2 = 437
$ws.Range("I6").Value2 = 258.92856
$ws.Range("J6").Value2 = 1268
$ws.Range("K6").Value2 = 776.78568
$ws.Range("L6").Value2 = 3804
$ws.Range("M6").Value2 = -664.78568
$ws.Range("N6").Value2 = -4028

# ALC row 17
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value2 = 1792.0566
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 1792.0566
$ws.Range("K17").Value2 = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value2 = 5376.1698
$ws.Range("N17").Value2 = -5712.1698

# ALC row 33
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value2 = 38462016
$ws.Range("I33").Value2 = 146.21053
$ws.Range("J33").Value2 = 142858530
$ws.Range("K33").Value2 = 146.21053
$ws.Range("L33").Value2 = 142858530
$ws.Range("M33").Value2 = 82.78946999999999
$ws.Range("N33").Value2 = -142858988

# ALC row 88
$ws = $wb.Worksheets.Item(1)
$ws.Range("I88").Value2 = 10001.5
$ws.Range("J88").Value2 = 7355.8887
$ws.Range("K88").Value2 = 10001.5
$ws.Range("L88").Value2 = 7355.8887
$ws.Range("M88").Value2 = -9595.5
$ws.Range("N88").Value2 = -8167.8887

# ALC row 91
$ws = $wb.Worksheets.Item(1)
$ws.Range("I91").Value2 = 10001.5
$ws.Range("J91").Value2 = 7355.8887
$ws.Range("K91").Value2 = 10001.5
$ws.Range("L91").Value2 = 7355.8887
$ws.Range("M91").Value2 = -8597.5
$ws.Range("N91").Value2 = -10163.8887

# ALC row 110
$ws = $wb.Worksheets.Item(1)
$ws.Range("H110").Value2 = 37459.1
$ws.Range("J110").Value2 = 37459.1
$ws.Range("L110").Value2 = 37459.1
$ws.Range("N110").Value2 = -45639.1

# ALC row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value2 = 3935.7334
$ws.Range("I132").Value2 = 1855.0968
$ws.Range("K132").Value2 = 5565.2904
$ws.Range("M132").Value2 = -3035.2904

# ALC row 137
$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value2 = 1425852.4
$ws.Range("I137").Value2 = 1200.2354
$ws.Range("K137").Value2 = 3600.7062
$ws.Range("M137").Value2 = -1050.7062

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value2 = 2395730.8
$ws.Range("I138").Value2 = 1048.1111
$ws.Range("J138").Value2 = 5158826
$ws.Range("K138").Value2 = 3144.3333
$ws.Range("L138").Value2 = 15476478
$ws.Range("M138").Value2 = 1995.6667
$ws.Range("N138").Value2 = -15486758

# ARM row 8
$ws = $wb.Worksheets.Item(2)
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 0
$ws.Range("K8").Value2 = 0
$ws.Range("M8").ClearContents()

# BSM row 11
$ws = $wb.Worksheets.Item(3)
$ws.Range("H11").Value2 = 266.66666
$ws.Range("I11").Value2 = 266.66666
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 266.66666
$ws.Range("L11").Value2 = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value2 = -126.66666

# BSM row 86
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value2 = 1512.9678
$ws.Range("J86").Value2 = 1637.6666
$ws.Range("L86").Value2 = 1637.6666
$ws.Range("N86").Value2 = -3883.6666

# BSM row 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value2 = 1512.9678
$ws.Range("J89").Value2 = 1637.6666
$ws.Range("L89").Value2 = 8188.333000000001
$ws.Range("N89").Value2 = -19420.333

# CRP row 22
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value2 = 736.25
$ws.Range("I22").Value2 = 241.42857
$ws.Range("J22").Value2 = 1002.6923
$ws.Range("K22").Value2 = 241.42857
$ws.Range("L22").Value2 = 1002.6923
$ws.Range("M22").Value2 = 108.57143
$ws.Range("N22").Value2 = -1702.6923

# CRP row 58
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value2 = 4503.2334
$ws.Range("I58").Value2 = 4744.3335
$ws.Range("J58").Value2 = 2333.3333
$ws.Range("K58").Value2 = 4744.3335
$ws.Range("L58").Value2 = 2333.3333
$ws.Range("M58").Value2 = -4541.3335
$ws.Range("N58").Value2 = -2739.3333

# CRP row 94
$ws = $wb.Worksheets.Item(4)
$ws.Range("H94").Value2 = 1693.4667
$ws.Range("I94").Value2 = 1698.7142
$ws.Range("J94").Value2 = 1688.875
$ws.Range("K94").Value2 = 1698.7142
$ws.Range("L94").Value2 = 1688.875
$ws.Range("M94").Value2 = -1247.7142
$ws.Range("N94").Value2 = -2590.875

# CRP row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value2 = 1096.8846
$ws.Range("I134").Value2 = 1085.9318
$ws.Range("J134").Value2 = 1157.125
$ws.Range("K134").Value2 = 3257.7954
$ws.Range("L134").Value2 = 3471.375
$ws.Range("M134").Value2 = -722.7954
$ws.Range("N134").Value2 = -8541.375

# CRP row 136
$ws = $wb.Worksheets.Item(4)
$ws.Range("H136").Value2 = 4503.2334
$ws.Range("I136").Value2 = 4744.3335
$ws.Range("J136").Value2 = 2333.3333
$ws.Range("K136").Value2 = 14233.0005
$ws.Range("L136").Value2 = 6999.999899999999
$ws.Range("M136").Value2 = -11683.0005
$ws.Range("N136").Value2 = -12099.9999

# CUL row 4
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value2 = 935.7368
$ws.Range("I4").Value2 = 319.6
$ws.Range("J4").Value2 = 1155.7858
$ws.Range("K4").Value2 = 958.8000000000001
$ws.Range("L4").Value2 = 3467.3574
$ws.Range("M4").Value2 = -846.8000000000001
$ws.Range("N4").Value2 = -3691.3574

# CUL row 6
$ws = $wb.Worksheets.Item(5)
$ws.Range("H6").Value2 = 877.94116
$ws.Range("I6").Value2 = 772.9286
$ws.Range("J6").Value2 = 1368
$ws.Range("K6").Value2 = 2318.7858
$ws.Range("L6").Value2 = 4104
$ws.Range("M6").Value2 = -2205.7858
$ws.Range("N6").Value2 = -4330

# CUL row 12
$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value2 = 272.02942
$ws.Range("I12").Value2 = 388.58334
$ws.Range("J12").Value2 = 208.45454
$ws.Range("K12").Value2 = 1165.75002
$ws.Range("L12").Value2 = 625.3636200000001
$ws.Range("M12").Value2 = -992.7500199999999
$ws.Range("N12").Value2 = -971.3636200000001

# CUL row 17
$ws = $wb.Worksheets.Item(5)
$ws.Range("H17").Value2 = 475.25
$ws.Range("I17").Value2 = 349.5
$ws.Range("J17").Value2 = 601
$ws.Range("K17").Value2 = 1048.5
$ws.Range("L17").Value2 = 1803
$ws.Range("M17").Value2 = -879.5
$ws.Range("N17").Value2 = -2141

# CUL row 33
$ws = $wb.Worksheets.Item(5)
$ws.Range("H33").Value2 = 63.57143
$ws.Range("I33").Value2 = 47.75
$ws.Range("J33").Value2 = 84.666664
$ws.Range("K33").Value2 = 286.5
$ws.Range("L33").Value2 = 507.999984
$ws.Range("M33").Value2 = -3.5
$ws.Range("N33").Value2 = -1073.999984

# CUL row 46
$ws = $wb.Worksheets.Item(5)
$ws.Range("H46").Value2 = 751808.4399999999
$ws.Range("I46").Value2 = 569.8570999999999
$ws.Range("J46").Value2 = 1156321.5
$ws.Range("K46").Value2 = 1709.5713
$ws.Range("L46").Value2 = 3468964.5
$ws.Range("M46").Value2 = -1618.5713
$ws.Range("N46").Value2 = -3469146.5

# CUL row 61
$ws = $wb.Worksheets.Item(5)
$ws.Range("H61").Value2 = 474.35294
$ws.Range("I61").Value2 = 275.66666
$ws.Range("J61").Value2 = 582.7273
$ws.Range("K61").Value2 = 826.9999799999999
$ws.Range("L61").Value2 = 1748.1819
$ws.Range("M61").Value2 = -611.9999799999999
$ws.Range("N61").Value2 = -2178.1819

# CUL row 70
$ws = $wb.Worksheets.Item(5)
$ws.Range("H70").Value2 = 4273.8184
$ws.Range("I70").Value2 = 1836
$ws.Range("J70").Value2 = 5961.5386
$ws.Range("K70").Value2 = 5508
$ws.Range("L70").Value2 = 17884.6158
$ws.Range("M70").Value2 = -5193
$ws.Range("N70").Value2 = -18514.6158

# CUL row 73
$ws = $wb.Worksheets.Item(5)
$ws.Range("H73").Value2 = 4273.8184
$ws.Range("I73").Value2 = 1836
$ws.Range("J73").Value2 = 5961.5386
$ws.Range("K73").Value2 = 5508
$ws.Range("L73").Value2 = 17884.6158
$ws.Range("M73").Value2 = -4416
$ws.Range("N73").Value2 = -20068.6158

# CUL row 80
$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value2 = 1654.8
$ws.Range("I80").Value2 = 654.4286
$ws.Range("J80").Value2 = 2193.4614
$ws.Range("K80").Value2 = 1963.2858
$ws.Range("L80").Value2 = 6580.3842
$ws.Range("M80").Value2 = -1027.2858
$ws.Range("N80").Value2 = -8452.3842

# CUL row 83
$ws = $wb.Worksheets.Item(5)
$ws.Range("H83").Value2 = 1654.8
$ws.Range("I83").Value2 = 654.4286
$ws.Range("J83").Value2 = 2193.4614
$ws.Range("K83").Value2 = 5889.8574
$ws.Range("L83").Value2 = 19741.1526
$ws.Range("M83").Value2 = -1209.8574
$ws.Range("N83").Value2 = -29101.1526

# CUL row 113
$ws = $wb.Worksheets.Item(5)
$ws.Range("H113").Value2 = 777424.9399999999
$ws.Range("I113").Value2 = 1782975.5
$ws.Range("J113").Value2 = 408.54544
$ws.Range("K113").Value2 = 5348926.5
$ws.Range("L113").Value2 = 1225.63632
$ws.Range("M113").Value2 = -5346756.5
$ws.Range("N113").Value2 = -5565.63632

# CUL row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value2 = 1000.19403
$ws.Range("I131").Value2 = 0
$ws.Range("J131").Value2 = 1000.19403
$ws.Range("K131").Value2 = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value2 = 3000.58209
$ws.Range("N131").Value2 = -13080.58209

# CUL row 138
$ws = $wb.Worksheets.Item(5)
$ws.Range("H138").Value2 = 1402.4
$ws.Range("I138").Value2 = 1072.5
$ws.Range("J138").Value2 = 2722
$ws.Range("K138").Value2 = 3217.5
$ws.Range("L138").Value2 = 8166
$ws.Range("M138").Value2 = 1922.5
$ws.Range("N138").Value2 = -18446

# GSM row 4
$ws = $wb.Worksheets.Item(6)
$ws.Range("H4").Value2 = 823.0833
$ws.Range("I4").Value2 = 1293.3334
$ws.Range("J4").Value2 = 666.3333
$ws.Range("K4").Value2 = 1293.3334
$ws.Range("L4").Value2 = 666.3333
$ws.Range("M4").Value2 = -1181.3334
$ws.Range("N4").Value2 = -890.3333

# GSM row 80
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value2 = 2318
$ws.Range("I80").Value2 = 2343.5293
$ws.Range("J80").Value2 = 2101
$ws.Range("K80").Value2 = 2343.5293
$ws.Range("L80").Value2 = 2101
$ws.Range("M80").Value2 = -1345.5293
$ws.Range("N80").Value2 = -4097

# GSM row 83
$ws = $wb.Worksheets.Item(6)
$ws.Range("H83").Value2 = 2318
$ws.Range("I83").Value2 = 2343.5293
$ws.Range("J83").Value2 = 2101
$ws.Range("K83").Value2 = 11717.6465
$ws.Range("L83").Value2 = 10505
$ws.Range("M83").Value2 = -6725.646500000001
$ws.Range("N83").Value2 = -20489

# LTW row 69
$ws = $wb.Worksheets.Item(7)
$ws.Range("H69").Value2 = 29537
$ws.Range("I69").Value2 = 30148
$ws.Range("J69").Value2 = 29333.334
$ws.Range("K69").Value2 = 30148
$ws.Range("L69").Value2 = 29333.334
$ws.Range("M69").Value2 = -29337
$ws.Range("N69").Value2 = -30955.334

# LTW row 72
$ws = $wb.Worksheets.Item(7)
$ws.Range("H72").Value2 = 29537
$ws.Range("I72").Value2 = 30148
$ws.Range("J72").Value2 = 29333.334
$ws.Range("K72").Value2 = 90444
$ws.Range("L72").Value2 = 88000.00199999999
$ws.Range("M72").Value2 = -86388
$ws.Range("N72").Value2 = -96112.00199999999

# WVR row 10
$ws = $wb.Worksheets.Item(8)
$ws.Range("H10").Value2 = 10000
$ws.Range("I10").Value2 = 10000
$ws.Range("K10").Value2 = 10000
$ws.Range("M10").Value2 = -9831

# WVR row 81
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value2 = 6200
$ws.Range("I81").Value2 = 2400
$ws.Range("J81").Value2 = 10000
$ws.Range("K81").Value2 = 4800
$ws.Range("L81").Value2 = 20000
$ws.Range("M81").Value2 = -3739
$ws.Range("N81").Value2 = -22122

# WVR row 84
$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value2 = 6200
$ws.Range("I84").Value2 = 2400
$ws.Range("J84").Value2 = 10000
$ws.Range("K84").Value2 = 24000
$ws.Range("L84").Value2 = 100000
$ws.Range("M84").Value2 = -18696
$ws.Range("N84").Value2 = -110608

# WVR row 113
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value2 = 715.1579
$ws.Range("I113").Value2 = 654.3
$ws.Range("J113").Value2 = 782.7778
$ws.Range("K113").Value2 = 1962.9
$ws.Range("L113").Value2 = 2348.3334
$ws.Range("M113").Value2 = 207.1000000000001
$ws.Range("N113").Value2 = -6688.3334

# WVR row 122
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value2 = 87007.664
$ws.Range("I122").Value2 = 8509.166999999999
$ws.Range("J122").Value2 = 401001.66
$ws.Range("K122").Value2 = 25527.501
$ws.Range("L122").Value2 = 1203004.98
$ws.Range("M122").Value2 = -23077.501
$ws.Range("N122").Value2 = -1207904.98

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value2 = 3213.875
$ws.Range("I132").Value2 = 3591.75
$ws.Range("J132").Value2 = 2647.0625
$ws.Range("K132").Value2 = 10775.25
$ws.Range("L132").Value2 = 7941.1875
$ws.Range("M132").Value2 = -8245.25
$ws.Range("N132").Value2 = -13001.1875
